$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column A of the new rows holds date-like text ("dd-mm-yyyy"); mark the
# cells as Text first so Excel's auto-detection doesn't turn them into
# real date serial numbers, then restore the default "Normal" style so
# the cells end up with no explicit style index, just like the rest of
# the sheet.
$ws.Range("A84:A85").NumberFormat = "@"

# --- Row 84: 04-08-2021 ---
$ws.Range("A84").Value = "04-08-2021"
$ws.Range("B84").Value = 900000
$ws.Range("C84").Value = 901400
$ws.Range("D84").Value = 900000
$ws.Range("E84").Value = 638600
$ws.Range("F84").Value = 261400
$ws.Range("G84").Value = 0.79

# --- Row 85: 05-08-2021 ---
$ws.Range("A85").Value = "05-08-2021"
$ws.Range("B85").Value = 900000
$ws.Range("C85").Value = 586000
$ws.Range("D85").Value = 450000
$ws.Range("E85").Value = 244000
$ws.Range("F85").Value = 206000
$ws.Range("G85").Value = 0.82

$ws.Range("A84:A85").Style = "Normal"
